$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 19, shifting existing rows 19:80 down to 20:81
$ws.Rows.Item(19).Insert()

# Fill the new row 19 with the new weekly price observation
$ws.Range("A19").Value = 8
$ws.Range("B19").Value = "Terminal La Palmera de La Serena"
$ws.Range("C19").Value = "Coquimbo"
$ws.Range("D19").Value = 44659
$ws.Range("E19").Value = 4
$ws.Range("F19").Value = 100112052
$ws.Range("G19").Value = "Albahaca"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 800
$ws.Range("K19").Value = 5000
$ws.Range("L19").Value = 5500
$ws.Range("M19").Value = 5250
$ws.Range("N19").Value = "$/docena de matas"
$ws.Range("O19").Value = "Provincia del Elquí"
$ws.Range("P19").Value = 875
$ws.Range("Q19").Value = 6
$ws.Range("R19").Value = "Hortaliza"
